$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextCell "D2" '62.807.20'
Set-TextCell "E2" '  +2.93%  '

# Row 3: Ethereum
Set-TextCell "D3" '2.974.98'
Set-TextCell "E3" '  +2.12%  '

# Row 4: TetherUSD
Set-TextCell "D4" '0.999'
Set-TextCell "E4" '  -0.09%  '

# Row 5: BNB
Set-TextCell "D5" '596.61'
Set-TextCell "E5" '  +1.20%  '

# Row 6: Solana
Set-TextCell "D6" '144.71'
Set-TextCell "E6" '  -0.13%  '

# Row 7: USDC
Set-TextCell "D7" '1.00'
Set-TextCell "E7" '  -0.04%  '

# Row 8: LidoStakedEther
Set-TextCell "D8" '2.974.79'
Set-TextCell "E8" '  +2.17%  '

# Row 9: XRP
Set-TextCell "E9" '  -0.24%  '

# Row 10: Toncoin
Set-TextCell "D10" '7.32'
Set-TextCell "E10" '  +6.25%  '

# Row 11: Dogecoin
Set-TextCell "D11" '0.145'
Set-TextCell "E11" '  +2.45%  '

# Row 12: Cardano
Set-TextCell "D12" '0.448'
Set-TextCell "E12" '  +1.98%  '

# Row 13: ShibaInu
Set-TextCell "D13" '0.0000238'
Set-TextCell "E13" '  +5.67%  '

# Row 14: Avalanche
Set-TextCell "D14" '33.63'
Set-TextCell "E14" '  +0.52%  '

# Row 16: WrappedliquidstakedEther2.0
Set-TextCell "D16" '3.464.12'
Set-TextCell "E16" '  +2.01%  '

# Row 17: WrappedBTC
Set-TextCell "D17" '62.496.33'
Set-TextCell "E17" '  +2.61%  '

# Row 18: Polkadot
Set-TextCell "D18" '6.75'
Set-TextCell "E18" '  +0.93%  '

# Row 19: WrappedEther
Set-TextCell "D19" '2.970.10'
Set-TextCell "E19" '  +2.01%  '

# Row 20: BitcoinCash
Set-TextCell "D20" '443.25'
Set-TextCell "E20" '  +2.43%  '

# Row 21: Chainlink
Set-TextCell "D21" '13.63'
Set-TextCell "E21" '  +1.98%  '

# Row 22: Polygon
Set-TextCell "D22" '0.677'
Set-TextCell "E22" '  +0.14%  '

# Row 23: Uniswap
Set-TextCell "D23" '7.20'
Set-TextCell "E23" '  +1.29%  '

# Row 24: Litecoin
Set-TextCell "D24" '82.07'
Set-TextCell "E24" '  +0.73%  '

# Row 25: RenderToken
Set-TextCell "D25" '10.90'
Set-TextCell "E25" '  +0.64%  '

# Row 26: InternetComputer(DFINITY)
Set-TextCell "D26" '12.07'
Set-TextCell "E26" '  +2.28%  '

# Row 27: Fetch.AI
Set-TextCell "D27" '2.16'
Set-TextCell "E27" '  -1.77%  '

# Row 28: Dai
Set-TextCell "E28" '  -0.03%  '

# Row 29: PancakeSwap
Set-TextCell "D29" '2.62'
Set-TextCell "E29" '  +0.92%  '

# Row 30: NEARProtocol
Set-TextCell "D30" '7.04'
Set-TextCell "E30" '  +0.91%  '

# Row 31: ImmutableX
Set-TextCell "D31" '2.13'
Set-TextCell "E31" '  -6.51%  '

# Row 32: EthereumClassic
Set-TextCell "D32" '26.58'
Set-TextCell "E32" '  +0.17%  '

# Row 33: Hedera
Set-TextCell "D33" '0.108'
Set-TextCell "E33" '  -0.80%  '

# Row 34: FirstDigitalUSD
Set-TextCell "D34" '0.999'
Set-TextCell "E34" '  -0.10%  '

# Row 35: PEPE
Set-TextCell "D35" '0.0₃0884'
Set-TextCell "E35" '  +1.73%  '

# Row 36: Mantle
Set-TextCell "E36" '  -1.04%  '

# Row 37: Filecoin
Set-TextCell "D37" '5.65'
Set-TextCell "E37" '  +0.55%  '

# Row 38: OKB
Set-TextCell "E38" '  +0.36%  '

# Row 39: dogwifhat
Set-TextCell "E39" '  -1.10%  '

# Row 40: Stacks
Set-TextCell "D40" '2.01'
Set-TextCell "E40" '  +1.96%  '

# Row 41: Cosmos
Set-TextCell "E41" '  +0.99%  '

# Row 42: Kaspa
Set-TextCell "E42" '  -2.45%  '

# Row 43: TheGraph
Set-TextCell "D43" '0.282'
Set-TextCell "E43" '  -0.29%  '

# Row 44: Arweave
Set-TextCell "D44" '39.12'
Set-TextCell "E44" '  -4.42%  '

# Row 45: Bittensor
Set-TextCell "B45" 'Bittensor'
Set-TextCell "C45" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D45" '370.38'
Set-TextCell "E45" '  -2.06%  '

# Row 46: Maker
Set-TextCell "B46" 'Maker'
Set-TextCell "C46" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D46" '2.700.76'
Set-TextCell "E46" '  +0.12%  '

# Row 47: VeChain
Set-TextCell "D47" '0.0342'
Set-TextCell "E47" '  -1.46%  '

# Row 48: Monero
Set-TextCell "D48" '133.87'
Set-TextCell "E48" '  +0.20%  '

# Row 49: USDe
Set-TextCell "E49" '  +0.06%  '

# Row 50: InjectiveProtocol
Set-TextCell "E50" '  -2.20%  '

# Row 51: Stellar
Set-TextCell "E51" '  -0.65%  '
